$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw")

# ---------------------------------------------------------------------------
# Add a new optional attribute "d_max" / "float64" to the Plant table
# (new row 26, right after the existing "storage_capacity" entry).
# Written first so new shared-string entries come out in the same order
# as the authored workbook (d_max, inflows, inflow).
# ---------------------------------------------------------------------------
$ws.Range("G26").Value = "d_max"
$ws.Range("H26").Value = "float64"

# ---------------------------------------------------------------------------
# Add a new "inflows" attribute table in columns AE:AF (rows 1-8), mirroring
# the layout of the other small tables on the sheet (e.g. "ntc" in AB:AC).
# ---------------------------------------------------------------------------

# Header row (merged title cell)
$ws.Range("AE1").Value = "inflows"

# "attributes / type" sub-header
$ws.Range("AE2").Value = "attributes "
$ws.Range("AF2").Value = "type"

# data rows: index / any, timestep / any, plant / plants.index, inflow / float64
$ws.Range("AE3").Value = "index"
$ws.Range("AF3").Value = "any"

$ws.Range("AE4").Value = "timestep"
$ws.Range("AF4").Value = "any"

$ws.Range("AE5").Value = "plant"
$ws.Range("AF5").Value = "plants.index"

$ws.Range("AE6").Value = "inflow"
$ws.Range("AF6").Value = "float64"

# row 7 stays empty under this block (cells removed entirely)
$ws.Range("AE7:AF7").ClearContents() | Out-Null

# "optional attributes" sub-header (no optional attributes follow)
$ws.Range("AE8").Value = "optional attributes"

# The second cell of a title/header pair is left as an (already-shared)
# empty string, matching the rest of the sheet's other blocks.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("AF1").PasteSpecial(-4163) | Out-Null   # xlPasteValues

$ws.Range("W8").Copy() | Out-Null
$ws.Range("AF8").PasteSpecial(-4163) | Out-Null   # xlPasteValues

$ws.Range("AE1:AF1").Merge() | Out-Null

# Copy formatting (fill/border/alignment) from the analogous cells of other
# tables so the new block looks like the rest of the sheet.
$ws.Range("G1:H1").Copy() | Out-Null
$ws.Range("AE1:AF1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("G2:H2").Copy() | Out-Null
$ws.Range("AE2:AF2").PasteSpecial(-4122) | Out-Null

$ws.Range("G16:H16").Copy() | Out-Null
$ws.Range("AE8:AF8").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Misc view state that changed in the authored workbook.
# ---------------------------------------------------------------------------
$ws.Range("AB14").Select() | Out-Null
